$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read current (pre-edit) values for the columns that differ across rows 2-4 ---
# Column A (Id, numeric), I (Antal, stored as text), K (Alder-Stadium, text),
# Q (Ost, numeric), R (Nord, numeric), AC (Publik kommentar, text)

$A2 = $ws.Range("A2").Value()
$I2 = $ws.Range("I2").Value()
$K2 = $ws.Range("K2").Value()
$Q2 = $ws.Range("Q2").Value()
$R2 = $ws.Range("R2").Value()
$AC2 = $ws.Range("AC2").Value()

$A3 = $ws.Range("A3").Value()
$I3 = $ws.Range("I3").Value()
$K3 = $ws.Range("K3").Value()
$Q3 = $ws.Range("Q3").Value()
$R3 = $ws.Range("R3").Value()
$AC3 = $ws.Range("AC3").Value()

$A4 = $ws.Range("A4").Value()
$I4 = $ws.Range("I4").Value()
$K4 = $ws.Range("K4").Value()
$Q4 = $ws.Range("Q4").Value()
$R4 = $ws.Range("R4").Value()
$AC4 = $ws.Range("AC4").Value()

# --- Make sure column I keeps a text number format, so numeric-looking
# strings ("10", "25", "500") are not silently re-typed as numbers.
# Column K never holds a numeric-looking value, so it needs no such hint. ---
$ws.Range("I2:I4").NumberFormat = "@"

# --- Write the rotated values: row2 -> row3, row3 -> row4, row4 -> row2 ---

# New row 2 gets old row 4's values
$ws.Range("A2").Value = $A4
$ws.Range("I2").Value = $I4
$ws.Range("K2").Value = $K4
$ws.Range("Q2").Value = $Q4
$ws.Range("R2").Value = $R4
$ws.Range("AC2").Value = $AC4

# New row 3 gets old row 2's values
$ws.Range("A3").Value = $A2
$ws.Range("I3").Value = $I2
$ws.Range("K3").Value = $K2
$ws.Range("Q3").Value = $Q2
$ws.Range("R3").Value = $R2
$ws.Range("AC3").Value = $AC2

# New row 4 gets old row 3's values
$ws.Range("A4").Value = $A3
$ws.Range("I4").Value = $I3
$ws.Range("K4").Value = $K3
$ws.Range("Q4").Value = $Q3
$ws.Range("R4").Value = $R3
$ws.Range("AC4").Value = $AC3
